$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Separator row (row 9): copy the header row's format (fill +
#    alignment) down onto A9:C9, leaving the cells empty - this is a
#    blank colored spacer row between the Jan 28 and Jan 29 entries.
# ------------------------------------------------------------------
$ws.Range("A1:C1").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the theme fill color (no tint) to both the header row and
# the new separator row so they share the same, fully-saturated
# accent color.
$ws.Range("A1:C1").Interior.ThemeColor = 8
$ws.Range("A9:C9").Interior.ThemeColor = 8

# ------------------------------------------------------------------
# 2) Prime the style "templates" we will reuse for the new data rows
#    by copying formats from existing, equivalently-styled cells.
#    Column A / C (plain center/center) -> copy from A3 / C3
#    Column B, no-wrap (left/center)    -> copy from B3
#    Column B, wrap (left/center/wrap)  -> copy from B4
# ------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A10:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C3").Copy()
$ws.Range("C10:C19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B14:B18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B4").Copy()
$ws.Range("B11:B13").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Values for the new timesheet entries (29-Jan-2020).
#    The Timestamp/Task text is entered in the same order the author
#    actually typed it (not strictly top-to-bottom - a couple of Task
#    notes were written out of order relative to their Timestamp row),
#    so that the shared-string table comes out in the same sequence.
# ------------------------------------------------------------------
$ws.Range("A10").Value = "Jan 29 10:00 to 11:00"
$ws.Range("A11").Value = "Jan 29 11:00 to 12:00"
$ws.Range("B11").Value = "Modified code of saving output graphs, now graphs are not saved.`nApplied two pair z-test"
$ws.Range("B10").Value = "Applied 1 paired z-test on good day and data file values. "

$ws.Range("A12").Value = "Jan 29 12:00 to 13:00"
$ws.Range("B12").Value = "Applied f-test. Applying anova test, modifieying data according to`nanova table. Performed missing values imputation."
$ws.Range("B13").Value = "Performed Anova test, printed anova table. Applied anova test of`nassumptions by using Levene test and Shapiro-Wilk test."
$ws.Range("A13").Value = "Jan 29 13:00 to 13:30"

$ws.Range("A14").Value = "Jan 29 13:30 to 14:00"
$ws.Range("B14").Value = "Observation and result documentation of statistical tests."

$ws.Range("A15").Value = "Jan 29 14:00 to 15:00"
$ws.Range("B15").Value = "Python Class"

$ws.Range("A16").Value = "Jan 29 15:00 to 16:00"
$ws.Range("B16").Value = "Searching slution for doing timeseries analysis"

$ws.Range("A17").Value = "Jan 29 16:00 to 17:00"
$ws.Range("B17").Value = "Visualization of both machine timeseries"

$ws.Range("A18").Value = "Jan 29 17:00 to 18:00"
$ws.Range("B18").Value = "Timeseries decomposition, modeified code of timeseries saving."

$ws.Range("A19").Value = "Jan 29 18:00 to 19:00"
$ws.Range("B19").Value = "Found errors when decomposing timeseries data. Tryed converting`ndata by interpolation and converting date to datetime."

# Location column - same "Infimetrics" shared string already exists
# in the workbook, so write order here is not significant.
$ws.Range("C10").Value = "Infimetrics"
$ws.Range("C11").Value = "Infimetrics"
$ws.Range("C12").Value = "Infimetrics"
$ws.Range("C13").Value = "Infimetrics"
$ws.Range("C14").Value = "Infimetrics"
$ws.Range("C15").Value = "Infimetrics"
$ws.Range("C16").Value = "Infimetrics"
$ws.Range("C17").Value = "Infimetrics"
$ws.Range("C18").Value = "Infimetrics"
$ws.Range("C19").Value = "Infimetrics"

# ------------------------------------------------------------------
# 4) Row heights for the wrapped, multi-line entries
# ------------------------------------------------------------------
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 45

# ------------------------------------------------------------------
# 5) Restore the column B style index/width bookkeeping and update
#    the view selection to match where the author left off editing.
# ------------------------------------------------------------------
$ws.Range("F17").Select()
